# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to Sheets/Aegis_Profits.xlsx
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR workbook tabs)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 93132
$ws.Range("I6").Value = 246800.33
$ws.Range("K6").Value = 740400.99
$ws.Range("M6").Value = -740288.99

$ws.Range("H9").Value = 315.16666
$ws.Range("I9").Value = 374.2
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 374.2
$ws.Range("L9").Value = 20
$ws.Range("M9").Value = -205.2
$ws.Range("N9").Value = -358

$ws.Range("H12").Value = 125125
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 166800
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 166800
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = -167140

$ws.Range("H21").Value = 11990
$ws.Range("I21").Value = 11990
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 11990
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -11522
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 11990
$ws.Range("I23").Value = 11990
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 11990
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -11756
$ws.Range("N23").ClearContents()

$ws.Range("H29").Value = 1290
$ws.Range("I29").Value = 1290
$ws.Range("K29").Value = 3870
$ws.Range("M29").Value = -3589

$ws.Range("H38").Value = 1897805.4
$ws.Range("I38").Value = 2688259.2
$ws.Range("J38").Value = 716
$ws.Range("K38").Value = 8064777.600000001
$ws.Range("L38").Value = 2148
$ws.Range("M38").Value = -8064405.600000001
$ws.Range("N38").Value = -2892

$ws.Range("H58").Value = 2452336.5
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9300

$ws.Range("H64").Value = 69666.13
$ws.Range("J64").Value = 3427.8572
$ws.Range("L64").Value = 3427.8572
$ws.Range("N64").Value = -3923.8572

$ws.Range("H67").Value = 69666.13
$ws.Range("J67").Value = 3427.8572
$ws.Range("L67").Value = 3427.8572
$ws.Range("N67").Value = -5143.8572

$ws.Range("H76").Value = 5317.5454
$ws.Range("I76").Value = 4831.5
$ws.Range("K76").Value = 4831.5
$ws.Range("M76").Value = -4516.5

$ws.Range("H79").Value = 5317.5454
$ws.Range("I79").Value = 4831.5
$ws.Range("K79").Value = 4831.5
$ws.Range("M79").Value = -3739.5

$ws.Range("H87").Value = 31878.5
$ws.Range("J87").Value = 31878.5
$ws.Range("L87").Value = 31878.5
$ws.Range("N87").Value = -34374.5

$ws.Range("H90").Value = 31878.5
$ws.Range("J90").Value = 31878.5
$ws.Range("L90").Value = 95635.5
$ws.Range("N90").Value = -108115.5

$ws.Range("H135").Value = 4049.121
$ws.Range("I135").Value = 1109.6666
$ws.Range("J135").Value = 5728.8096
$ws.Range("K135").Value = 9986.999400000001
$ws.Range("L135").Value = 51559.2864
$ws.Range("M135").Value = -7451.999400000001
$ws.Range("N135").Value = -56629.2864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28545.29
$ws.Range("I32").Value = 10026.06
$ws.Range("J32").Value = 106094.56
$ws.Range("K32").Value = 10026.06
$ws.Range("L32").Value = 106094.56
$ws.Range("M32").Value = -9739.059999999999
$ws.Range("N32").Value = -106668.56

$ws.Range("H61").Value = 2118.138
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 2131.7693
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 2131.7693
$ws.Range("M61").Value = -1788
$ws.Range("N61").Value = -2555.7693

$ws.Range("H74").Value = 1289.931
$ws.Range("I74").Value = 1288.4
$ws.Range("K74").Value = 1288.4
$ws.Range("M74").Value = -414.4000000000001

$ws.Range("H77").Value = 1289.931
$ws.Range("I77").Value = 1288.4
$ws.Range("K77").Value = 6442
$ws.Range("M77").Value = -2074

$ws.Range("H136").Value = 2118.138
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 2131.7693
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 6395.3079
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -11495.3079

$ws.Range("H138").Value = 57950
$ws.Range("J138").Value = 57950
$ws.Range("L138").Value = 57950
$ws.Range("N138").Value = -68230

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2700
$ws.Range("J62").Value = 2700
$ws.Range("L62").Value = 2700
$ws.Range("N62").Value = -3948

$ws.Range("H65").Value = 2700
$ws.Range("J65").Value = 2700
$ws.Range("L65").Value = 13500
$ws.Range("N65").Value = -19740

$ws.Range("H68").Value = 18060.143
$ws.Range("J68").Value = 18060.143
$ws.Range("L68").Value = 18060.143
$ws.Range("N68").Value = -19558.143

$ws.Range("H71").Value = 18060.143
$ws.Range("J71").Value = 18060.143
$ws.Range("L71").Value = 54180.429
$ws.Range("N71").Value = -61668.429

$ws.Range("H74").Value = 38177.5
$ws.Range("J74").Value = 38177.5
$ws.Range("L74").Value = 38177.5
$ws.Range("N74").Value = -39925.5

$ws.Range("H77").Value = 38177.5
$ws.Range("J77").Value = 38177.5
$ws.Range("L77").Value = 114532.5
$ws.Range("N77").Value = -123268.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1977.6923
$ws.Range("I136").Value = 1810.8334
$ws.Range("K136").Value = 5432.5002
$ws.Range("M136").Value = -332.5002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 107195.75
$ws.Range("I70").Value = 149279
$ws.Range("J70").Value = 9001.5
$ws.Range("K70").Value = 149279
$ws.Range("L70").Value = 9001.5
$ws.Range("M70").Value = -149009
$ws.Range("N70").Value = -9541.5

$ws.Range("H73").Value = 107195.75
$ws.Range("I73").Value = 149279
$ws.Range("J73").Value = 9001.5
$ws.Range("K73").Value = 149279
$ws.Range("L73").Value = 9001.5
$ws.Range("M73").Value = -148343
$ws.Range("N73").Value = -10873.5

$ws.Range("H80").Value = 125276800
$ws.Range("I80").Value = 200441500
$ws.Range("J80").Value = 2316.6667
$ws.Range("K80").Value = 200441500
$ws.Range("L80").Value = 2316.6667
$ws.Range("M80").Value = -200440502
$ws.Range("N80").Value = -4312.6667

$ws.Range("H83").Value = 125276800
$ws.Range("I83").Value = 200441500
$ws.Range("J83").Value = 2316.6667
$ws.Range("K83").Value = 1002207500
$ws.Range("L83").Value = 11583.3335
$ws.Range("M83").Value = -1002202508
$ws.Range("N83").Value = -21567.3335

$ws.Range("H132").Value = 3252.682
$ws.Range("I132").Value = 2268.353
$ws.Range("K132").Value = 6805.059
$ws.Range("M132").Value = -4275.059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 23262.223
$ws.Range("J123").Value = 23262.223
$ws.Range("L123").Value = 23262.223
$ws.Range("N123").Value = -33062.223

$ws.Range("H136").Value = 2324.8438
$ws.Range("I136").Value = 1690.9546
$ws.Range("K136").Value = 5072.8638
$ws.Range("M136").Value = -2522.8638
